$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Modelo" header in F1, reusing the bold/bordered/centered header
# style already used by A1:E1 (copy format from E1, then set the text)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"

# Update the recalculated metric values in row 2
$ws.Range("B2").Value = 0.2817704156383548
$ws.Range("C2").Value = 0.9945793637595082
$ws.Range("D2").Value = 0.4327639265559824

# Add the new model description in F2
$ws.Range("F2").Value = "Pipeline(steps=[('model', RandomForestRegressor(max_depth=5))])"
